# Update "想去人数" (F column) figures across the 展览 / 演出 / 全部类型 sheets
# to reflect the newly generated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2986
$ws.Range("F4").Value = 2540
$ws.Range("F7").Value = 78
$ws.Range("F9").Value = 2940
$ws.Range("F10").Value = 359
$ws.Range("F12").Value = 7565
$ws.Range("F13").Value = 361
$ws.Range("F15").Value = 109
$ws.Range("F16").Value = 256
$ws.Range("F19").Value = 9252
$ws.Range("F22").Value = 269
$ws.Range("F31").Value = 72
$ws.Range("F32").Value = 117
$ws.Range("F36").Value = 1417
$ws.Range("F37").Value = 1487
$ws.Range("F39").Value = 3938
$ws.Range("F41").Value = 44
$ws.Range("F42").Value = 1199
$ws.Range("F43").Value = 96
$ws.Range("F44").Value = 33
$ws.Range("F45").Value = 247
$ws.Range("F46").Value = 12
$ws.Range("F47").Value = 65
$ws.Range("F48").Value = 36
$ws.Range("F49").Value = 59

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 52
$ws.Range("F7").Value = 142
$ws.Range("F15").Value = 5
$ws.Range("F20").Value = 34

# --- Sheet: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 52
$ws.Range("F3").Value = 2986
$ws.Range("F7").Value = 2540
$ws.Range("F8").Value = 142
$ws.Range("F11").Value = 78
$ws.Range("F13").Value = 2940
$ws.Range("F14").Value = 359
$ws.Range("F18").Value = 7565
$ws.Range("F19").Value = 361
$ws.Range("F21").Value = 109
$ws.Range("F22").Value = 256
$ws.Range("F24").Value = 9252
$ws.Range("F25").Value = 269
$ws.Range("F32").Value = 72
$ws.Range("F33").Value = 117
$ws.Range("F36").Value = 1419
$ws.Range("F37").Value = 1487
$ws.Range("F40").Value = 3938
$ws.Range("F42").Value = 44
$ws.Range("F44").Value = 1199
$ws.Range("F45").Value = 96
$ws.Range("F46").Value = 247
$ws.Range("F47").Value = 65
$ws.Range("F48").Value = 36
$ws.Range("F49").Value = 59
